# Apply marksheet corrections: update "Marking" and "Total" rows
# B11: 3 -> 5 (Marking / Right)
# B12: 30 -> 50 (Total / Right)
# E12: "26/84" -> "50/140" (Total / Max, correct-over-total text)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 50
$ws.Range("E12").Value = "50/140"
